# Update template FAQ to clarify affiliations sheet
#
# In the "Users by Affiliations" block (rows 16-18), insert a new row for
# a "Month - Year" field and reword the existing "Totals" row definition
# to mention "new users" instead of just "users".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing "Totals" row (row 17) to make room
# for the new "Month - Year" field row.
$ws.Rows("17:17").Insert()

# Copy the formatting from the row above (Affiliation row, still intact)
# onto the freshly inserted blank row so the new row's cell styles match
# the rest of the section (A uses style 14, B uses style 6, C uses style 7).
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)

# Update the "Totals" row's definition text (now shifted down to row 18)
# to refer to "new users" rather than "users".
$ws.Range("C18").Value2 = "The number of new users by affiliation over the one year time period"

# Fill in the new "Month - Year" row (row 17).
$ws.Range("B17").Value2 = "Month - Year"
$ws.Range("C17").Value2 = "The number of new users by affiliation over the month"

# Match the saved selection shown in the target workbook.
$ws.Range("C18").Select()
